$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("threat_list")

$ws.Range("A27").Value = "TR-61"
$ws.Range("B27").Value = "Brainstorming"
$ws.Range("C27").Value = "Tampering/ Information Disclosure/ Spoofing"
$ws.Range("D27").Value = "Server <=> Client"
$ws.Range("E27").Value = "By changing the server/client's certificate or key, an attacker may attempt to connect to an unauthorized client.`nAnd attacker can try to steal the information of the encryption channel."
$ws.Range("F27").Value = "Need to protect or verify the certificates and keys used by the server and client"
$ws.Range("G27").Value = "A server and client program must perform an integrity check before using a certificate or key."
$ws.Range("G27").WrapText = $false
$ws.Range("H27").Value = "Certificate & Key file existance check`nIntegrity Check with hash function`n- Use OpenSSL library of latest version (1.1.1k)`n- Use an algorithm that are stronger than sha256"

$ws.Activate()
$ws.Range("G29").Select()
